$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.827.43"
$ws.Range("E2").Value = "  -0.81%  "
$ws.Range("D3").Value = "3.494.55"
$ws.Range("E3").Value = "  -1.99%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'602.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.89%  "
$ws.Range("D6").Value = "'197.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.09%  "
$ws.Range("D7").Value = "'0.623"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.70%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").Value = "'0.209"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.83%  "
$ws.Range("D10").Value = "'0.652"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.95%  "
$ws.Range("D11").Value = "'54.11"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("D12").Value = "'0.0000301"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.11%  "
$ws.Range("D13").Value = "'9.54"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.19%  "
$ws.Range("D14").Value = "4.054.36"
$ws.Range("E14").Value = "  -1.80%  "
$ws.Range("D15").Value = "'592.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.78%  "
$ws.Range("D16").Value = "69.897.28"
$ws.Range("E16").Value = "  -0.80%  "
$ws.Range("D17").Value = "'18.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.48%  "
$ws.Range("D18").Value = "'12.61"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.36%  "
$ws.Range("D19").Value = "3.498.17"
$ws.Range("E19").Value = "  -0.96%  "
$ws.Range("E20").Value = "  +0.15%  "
$ws.Range("D21").Value = "'0.984"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.25%  "
$ws.Range("D22").Value = "'17.90"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.06%  "
$ws.Range("D23").Value = "'103.71"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +9.90%  "
$ws.Range("D24").Value = "'4.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.72%  "
$ws.Range("D25").Value = "'4.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.15%  "
$ws.Range("D26").Value = "'3.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.36%  "
$ws.Range("D27").Value = "'10.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").Value = "'9.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.39%  "
$ws.Range("D29").Value = "'33.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.18%  "
$ws.Range("D30").Value = "'4.49"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +21.18%  "
$ws.Range("D31").Value = "'7.26"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.45%  "
$ws.Range("D32").Value = "'12.69"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.53%  "
$ws.Range("E33").Value = "  +1.06%  "
$ws.Range("D34").Value = "'63.53"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.31%  "
$ws.Range("D35").Value = "3.684.10"
$ws.Range("E35").Value = "  +4.42%  "
$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D36").Value = "'0.998"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.27%  "
$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").Value = "0.0₃0801"
$ws.Range("E37").Value = "  +1.59%  "
$ws.Range("D38").Value = "'514.70"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.11%  "
$ws.Range("D39").Value = "'0.390"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.81%  "
$ws.Range("D40").Value = "'2.97"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.46%  "
$ws.Range("D41").Value = "'36.53"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.59%  "
$ws.Range("D42").Value = "'3.52"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.01%  "
$ws.Range("D43").Value = "'0.137"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.39%  "
$ws.Range("D44").Value = "'0.0457"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.71%  "
$ws.Range("D45").Value = "'2.84"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.00%  "
$ws.Range("E46").Value = "  -0.98%  "
$ws.Range("E47").Value = "  -4.89%  "
$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D48").Value = "'1.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.58%  "
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").Value = "'8.75"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.14%  "
$ws.Range("D50").Value = "'132.23"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.41%  "
$ws.Range("D51").Value = "'0.000240"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.78%  "
